$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-02 Tuesday" "2024-07-03 Wednesday"

Replace-Text "281÷8=" "965÷3="
Replace-Text "689÷8=" "647÷9="
Replace-Text "676÷2=" "964÷5="
Replace-Text "746÷6=" "931÷8="
Replace-Text "998÷2=" "671÷6="

Replace-Text "893÷8=" "343÷3="
Replace-Text "150÷9=" "380÷9="
Replace-Text "863÷6=" "208÷8="
Replace-Text "992÷7=" "141÷9="
Replace-Text "973÷4=" "883÷6="

Replace-Text "731÷7=" "457÷2="
Replace-Text "197÷2=" "909÷7="
Replace-Text "311÷3=" "492÷9="
Replace-Text "428÷9=" "147÷2="
Replace-Text "884÷2=" "446÷7="

Replace-Text "523÷4=" "208÷6="
Replace-Text "843÷6=" "771÷4="
Replace-Text "516÷8=" "547÷7="
Replace-Text "319÷4=" "400÷3="
Replace-Text "247÷8=" "795÷2="

Replace-Text "134÷2=" "989÷7="
Replace-Text "685÷5=" "594÷9="
Replace-Text "220÷5=" "494÷5="
Replace-Text "600÷6=" "707÷7="
Replace-Text "615÷3=" "272÷4="
